# Updates diagnostics worksheet shared-string values per commit:
# "fix: Fix issue in fhir json generation for missing encounters #88"
#
# Column C (version) bumped from 0.10.1 -> 0.10.2 for all data rows (2-96),
# and column L (orch_session_issue_id) UUID values are regenerated for the
# affected issue rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C2" = "0.10.2"
    "L2" = "9fd3e417-3582-46ba-95c3-00879067e726"
    "C3" = "0.10.2"
    "L3" = "2dc6e40d-6554-43c6-aad8-c8a180ffd722"
    "C4" = "0.10.2"
    "L4" = "d1baf160-351e-4bc2-ba92-40e7d4170bd3"
    "C5" = "0.10.2"
    "L5" = "b7a39775-4cfe-417a-b8bb-206cccf78454"
    "C6" = "0.10.2"
    "L6" = "a3901164-c413-4cda-b0a4-c3ea8aef4832"
    "C7" = "0.10.2"
    "L7" = "303931bf-eb91-49b9-ad4f-d0e1beadb13a"
    "C8" = "0.10.2"
    "L8" = "e1d31992-a452-4822-b96a-918b282dc82d"
    "C9" = "0.10.2"
    "L9" = "6dce7fc7-a4e4-42e4-a971-a65dae0efc58"
    "C10" = "0.10.2"
    "L10" = "a75f49ac-025f-4b2a-a9cf-a1ed766f934e"
    "C11" = "0.10.2"
    "L11" = "60d973e9-5bd8-4e77-925b-0b7f61aabcb5"
    "C12" = "0.10.2"
    "L12" = "8529fa01-9219-4cd1-8a7a-7288a1b43a56"
    "C13" = "0.10.2"
    "L13" = "458691ba-4886-470a-86b4-40cb38d14df4"
    "C14" = "0.10.2"
    "L14" = "7e90c9a5-e4ad-41a1-8bd5-997e7ca346ad"
    "C15" = "0.10.2"
    "L15" = "fc045f6a-41d9-46db-a729-7f2af5061fcb"
    "C16" = "0.10.2"
    "L16" = "03775aee-13d6-4217-aee8-0cbd84ed2213"
    "C17" = "0.10.2"
    "L17" = "411db574-fcbc-4adb-a5b9-98304503b359"
    "C18" = "0.10.2"
    "L18" = "44904906-191b-4bd6-9f4e-e339985b4562"
    "C19" = "0.10.2"
    "L19" = "628e8f26-7094-4d02-b963-4a420f9980dc"
    "C20" = "0.10.2"
    "L20" = "7a16a168-aa27-426a-a2c8-bdedc085a3cc"
    "C21" = "0.10.2"
    "L21" = "f857ef5b-e9c2-440d-b0be-b5658daeb055"
    "C22" = "0.10.2"
    "L22" = "1c2a17e7-2c9e-4362-8bd9-2cb140510355"
    "C23" = "0.10.2"
    "L23" = "9a80c399-b3fa-4b35-b1a7-f4c311ad2d65"
    "C24" = "0.10.2"
    "L24" = "ffdaf977-7c72-480b-8a1f-486cf8c06938"
    "C25" = "0.10.2"
    "L25" = "40642f3e-d6b5-4d5a-8521-e2a2089881f8"
    "C26" = "0.10.2"
    "L26" = "fffdb8a8-65db-4f70-8a2a-c7cab63b1359"
    "C27" = "0.10.2"
    "L27" = "ca9da0c6-a511-437e-afd0-22d6746a845e"
    "C28" = "0.10.2"
    "L28" = "b9f32bba-e240-4125-8aad-170cafe91366"
    "C29" = "0.10.2"
    "L29" = "0a91b605-b54d-4eaa-ad16-08d2d4e86e49"
    "C30" = "0.10.2"
    "L30" = "76f680ab-b8d0-43d1-8566-346e7c636b3f"
    "C31" = "0.10.2"
    "L31" = "000b63ad-edf9-4db5-82c7-d75d4640c704"
    "C32" = "0.10.2"
    "L32" = "d2050cb1-dc13-41cd-8f8e-08759a473bac"
    "C33" = "0.10.2"
    "L33" = "0fe06e09-16ba-4983-8afc-7f5c3ba9878f"
    "C34" = "0.10.2"
    "L34" = "691c4eb1-c5a0-42fb-bf51-f713b271fbab"
    "C35" = "0.10.2"
    "L35" = "852d66ab-ebbe-4ceb-91f5-9c67749ada3a"
    "C36" = "0.10.2"
    "L36" = "cad7516f-cd2a-46c4-99d5-edb0bc26bbed"
    "C37" = "0.10.2"
    "L37" = "f99c8ba0-40cc-4561-9eb0-17ced4c9feb8"
    "C38" = "0.10.2"
    "L38" = "455a067d-70ab-495e-b295-c7bc2ee90ec8"
    "C39" = "0.10.2"
    "L39" = "9f48f3f4-fc97-4d57-bd3f-93bbc5324b12"
    "C40" = "0.10.2"
    "L40" = "43fc1838-ced3-4aa3-bb19-73098f3b9d62"
    "C41" = "0.10.2"
    "L41" = "5ae98409-8df3-4c15-9569-f03ee21b415c"
    "C42" = "0.10.2"
    "L42" = "4372558f-1383-4ecf-b839-e12bcce5a2b5"
    "C43" = "0.10.2"
    "L43" = "4bb0f8d1-23c6-4aca-8dbc-42e108b9fc06"
    "C44" = "0.10.2"
    "L44" = "5f445ad8-8db0-4925-b7bf-6c85df877131"
    "C45" = "0.10.2"
    "L45" = "bc53fd75-c097-4fc2-afe5-4e7cf0a3b064"
    "C46" = "0.10.2"
    "L46" = "3597a11a-0e5c-4d37-b6ad-f18b1b9ecbcf"
    "C47" = "0.10.2"
    "L47" = "09a9b03d-ccb8-42ce-8008-305a8734e3f5"
    "C48" = "0.10.2"
    "L48" = "5abf28d6-ae38-4903-9bf8-1d987127a05f"
    "C49" = "0.10.2"
    "L49" = "3f227d63-66bc-46e3-a857-df8762cd5aad"
    "C50" = "0.10.2"
    "L50" = "60825e9e-01d9-4a71-923d-7ff542de3859"
    "C51" = "0.10.2"
    "L51" = "0b3a2b49-665e-48fc-a339-8b6e36a946d7"
    "C52" = "0.10.2"
    "L52" = "efc679f5-9851-4c4d-98fb-79d4e5a53706"
    "C53" = "0.10.2"
    "L53" = "e146e43c-d072-435f-afd4-d257467eb7c5"
    "C54" = "0.10.2"
    "L54" = "e4283e08-dd82-4e8e-81f6-1459548407d3"
    "C55" = "0.10.2"
    "L55" = "f6566695-9564-412a-9b30-f1ca756ec324"
    "C56" = "0.10.2"
    "L56" = "e0046b1e-d7b7-4130-a750-12f5e9f6b1e9"
    "C57" = "0.10.2"
    "L57" = "b3630b75-695c-42fc-b6e3-4b1ed114c7cd"
    "C58" = "0.10.2"
    "L58" = "70118710-d1cf-451a-ad9f-b026a99f83f8"
    "C59" = "0.10.2"
    "L59" = "093771f8-b04f-4160-b211-96780736d3e9"
    "C60" = "0.10.2"
    "L60" = "24bc18d7-e8e8-4bdb-90db-c64184a559a3"
    "C61" = "0.10.2"
    "L61" = "4382f6da-5b1f-4ce2-950f-5090cbf9c8f3"
    "C62" = "0.10.2"
    "L62" = "0bc27811-4640-4391-ba21-aaa584255af4"
    "C63" = "0.10.2"
    "L63" = "c9d3df1f-f4f4-400b-ac90-42ab06c86eff"
    "C64" = "0.10.2"
    "L64" = "2fb74a0e-9bfa-4aa2-b077-bbc3ddc38483"
    "C65" = "0.10.2"
    "L65" = "b1c31bb6-6dd0-4ddc-a304-17e0648af64c"
    "C66" = "0.10.2"
    "L66" = "33617cda-dc6b-4010-a697-898da30b86d4"
    "C67" = "0.10.2"
    "L67" = "ed45fc53-dbeb-4e35-8ff9-89ddb40a05df"
    "C68" = "0.10.2"
    "L68" = "9c134df2-3342-4a86-a564-87f384badcde"
    "C69" = "0.10.2"
    "L69" = "fde877e1-3d43-48aa-880d-ee9d938409e1"
    "C70" = "0.10.2"
    "L70" = "8059d9d0-cea6-4ae2-b81c-4079faa11dbc"
    "C71" = "0.10.2"
    "L71" = "b02a7600-4e97-4088-af56-2caadfca079f"
    "C72" = "0.10.2"
    "L72" = "705e0968-d621-49f4-86d9-e7a495b978ea"
    "C73" = "0.10.2"
    "L73" = "1423673a-0da4-4ce7-a276-a2cda6102cdf"
    "C74" = "0.10.2"
    "L74" = "868799e5-46d1-4f7b-8f62-a661561d3aa5"
    "C75" = "0.10.2"
    "L75" = "4fad556f-7ba5-488d-b519-b2963704d522"
    "C76" = "0.10.2"
    "L76" = "23ce14b7-d7d1-4350-84e9-8951856defdf"
    "C77" = "0.10.2"
    "L77" = "fb9cc18f-b28b-45fc-99d1-040f368efb21"
    "C78" = "0.10.2"
    "L78" = "95aa3c19-e160-485e-a085-ac883cff8162"
    "C79" = "0.10.2"
    "L79" = "d8402606-7b3d-4892-809b-5bfae72697ca"
    "C80" = "0.10.2"
    "L80" = "a3f29ad4-58f4-4197-b8d9-47d5fd0032f6"
    "C81" = "0.10.2"
    "C82" = "0.10.2"
    "C83" = "0.10.2"
    "C84" = "0.10.2"
    "C85" = "0.10.2"
    "C86" = "0.10.2"
    "C87" = "0.10.2"
    "C88" = "0.10.2"
    "C89" = "0.10.2"
    "C90" = "0.10.2"
    "C91" = "0.10.2"
    "C92" = "0.10.2"
    "C93" = "0.10.2"
    "C94" = "0.10.2"
    "C95" = "0.10.2"
    "C96" = "0.10.2"
    "L96" = "33002cfe-91b5-4229-8808-b6204c335c5a"
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
